$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of data (2025-09-30) for both stations, matching the
# existing 四方坪站 / 高岭站 pattern used throughout the sheet.
$ws.Range("A60").Value = 45930
$ws.Range("B60").Value = "四方坪站"
$ws.Range("C60").Value = 10978.91
$ws.Range("D60").Value = 9163.65
$ws.Range("E60").Value = 3797.81
$ws.Range("F60").Value = 475

$ws.Range("A61").Value = 45930
$ws.Range("B61").Value = "高岭站"
$ws.Range("C61").Value = 5463.25
$ws.Range("D61").Value = 4403.68
$ws.Range("E61").Value = 1342.68
$ws.Range("F61").Value = 218

# Match number formats used by the rest of the columns (escaped exactly like
# the existing style definitions so Excel reuses the same style indices
# instead of minting new, duplicate ones).
$ws.Range("A60:A61").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("C60:E61").NumberFormat = "0.00_);[Red]\(0.00\)"
$ws.Range("F60:F61").NumberFormat = "0_);[Red]\(0\)"

# Update the view to reflect the scrolled/selected state after the edit.
$ws.Range("D64").Select()
